$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$values = @(
    "2025-12-19T07:01:23.444490+00:00",
    "2025-12-19T07:01:25.397645+00:00",
    "2025-12-19T07:01:25.397682+00:00",
    "2025-12-19T07:01:25.397700+00:00",
    "2025-12-19T07:01:27.792405+00:00",
    "2025-12-19T07:01:30.273346+00:00",
    "2025-12-19T07:01:32.667815+00:00",
    "2025-12-19T07:01:32.667844+00:00",
    "2025-12-19T07:01:34.597407+00:00",
    "2025-12-19T07:01:38.497990+00:00",
    "2025-12-19T07:01:38.498020+00:00",
    "2025-12-19T07:01:40.367563+00:00",
    "2025-12-19T07:01:42.259536+00:00",
    "2025-12-19T07:01:44.616872+00:00",
    "2025-12-19T07:01:46.970783+00:00",
    "2025-12-19T07:01:46.970811+00:00",
    "2025-12-19T07:01:46.970827+00:00",
    "2025-12-19T07:01:46.970843+00:00",
    "2025-12-19T07:01:48.888961+00:00",
    "2025-12-19T07:01:48.888988+00:00",
    "2025-12-19T07:01:48.889007+00:00",
    "2025-12-19T07:01:50.852003+00:00",
    "2025-12-19T07:01:50.852031+00:00",
    "2025-12-19T07:01:50.852048+00:00",
    "2025-12-19T07:01:53.178933+00:00",
    "2025-12-19T07:01:53.178961+00:00",
    "2025-12-19T07:01:55.090803+00:00",
    "2025-12-19T07:01:55.090834+00:00",
    "2025-12-19T07:01:55.090852+00:00",
    "2025-12-19T07:01:56.969942+00:00",
    "2025-12-19T07:01:58.854646+00:00",
    "2025-12-19T07:01:58.854676+00:00",
    "2025-12-19T07:02:03.727173+00:00",
    "2025-12-19T07:02:03.727203+00:00",
    "2025-12-19T07:02:05.589443+00:00",
    "2025-12-19T07:02:05.589470+00:00"
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 11).Value = $val
    $row = $row + 1
}
